# "Finished graph and report"
#
# The workbook's Sheet1 has Process Count (A), Execution Time in seconds (B)
# and a Normalized Execution Time formula column (C) that divides each B
# value by B2. Column B was left blank while the experiment was still
# running, so C currently evaluates to #DIV/0!. Now that the runs have
# finished, fill in the measured execution times - this lets column C (and
# the scatter chart plotted from A/C) calculate real values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$execTimes = @{
    2  = 2.5907710000000002
    3  = 1.308122
    4  = 0.68147599999999997
    5  = 0.34419100000000002
    6  = 0.260546
    7  = 0.253193
    8  = 0.241563
    9  = 0.24504200000000001
    10 = 0.24978400000000001
}

foreach ($row in $execTimes.Keys) {
    $ws.Cells.Item($row, 2).Value = $execTimes[$row]
}

# The chart was built off of some now-stale helper defined names
# (_xlchart.v1.*, _xlchart.v2.*) that Excel no longer needs once the chart
# itself carries its own cached series references - clear them out as part
# of tidying up the finished report.
$nameCount = $wb.Names.Count
for ($i = $nameCount; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}

# Leave the cursor where the author last clicked while reviewing the
# finished report.
$ws.Range("G6").Select() | Out-Null

$excel.CalculateFullRebuild() | Out-Null
